# Refresh the cryptos price/volume snapshot cells (D2:E51) with the
# latest scraped values, per the GitHub Actions update run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.493.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.87%  "
$ws.Range("D3").Value = "'2.413.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +8.77%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'323.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +11.74%  "
$ws.Range("D6").Value = "'104.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.88%  "
$ws.Range("D7").Value = "'0.645"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.62%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.659"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +11.52%  "
$ws.Range("D10").Value = "'42.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.99%  "
$ws.Range("D11").Value = "'0.0952"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.02%  "
$ws.Range("D12").Value = "'8.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("E13").Value = "  +5.09%  "
$ws.Range("D14").Value = "'17.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +17.02%  "
$ws.Range("E15").Value = "  +2.83%  "
$ws.Range("D16").Value = "'2.779.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.97%  "
$ws.Range("D17").Value = "'2.408.07"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.43%  "
$ws.Range("D18").Value = "'43.530.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.11%  "
$ws.Range("E19").Value = "  +6.01%  "
$ws.Range("E20").Value = "  +6.09%  "
$ws.Range("D21").Value = "'75.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.89%  "
$ws.Range("E22").Value = "  +3.94%  "
$ws.Range("D23").Value = "'261.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +14.18%  "
$ws.Range("E24").Value = "  +4.36%  "
$ws.Range("D25").Value = "'9.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.39%  "
$ws.Range("D26").Value = "'11.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.54%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "'3.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("D29").Value = "'22.82"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.27%  "
$ws.Range("D30").Value = "'179.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.72%  "
$ws.Range("E31").Value = "  +1.72%  "
$ws.Range("D32").Value = "'38.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.86%  "
$ws.Range("E33").Value = "  +3.11%  "
$ws.Range("D34").Value = "'0.0937"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.83%  "
$ws.Range("D35").Value = "'5.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.40%  "
$ws.Range("E36").Value = "  +6.43%  "
$ws.Range("D37").Value = "'4.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'0.0372"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("D39").Value = "'3.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.05%  "
$ws.Range("E40").Value = "  +3.86%  "
$ws.Range("E41").Value = "  +22.89%  "
$ws.Range("D42").Value = "'1.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +27.19%  "
$ws.Range("D43").Value = "'127.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +26.49%  "
$ws.Range("E44").Value = "  +2.91%  "
$ws.Range("D45").Value = "'69.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("D47").Value = "'12.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.07%  "
$ws.Range("D48").Value = "'9.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +14.96%  "
$ws.Range("D49").Value = "'5.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.50%  "
$ws.Range("E50").Value = "  +5.34%  "
$ws.Range("D51").Value = "'1.612.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.84%  "
